$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Elements"
# ---------------------------------------------------------------------------
$wsElements = $wb.Worksheets.Item("Elements")

# B8 value: "stock" -> "auxiliary"
$wsElements.Range("B8").Value = "auxiliary"

# Column widths
$wsElements.Columns.Item(1).ColumnWidth = 22.166666666666668
$wsElements.Columns.Item(2).ColumnWidth = 14.166666666666666

# Row heights for the data rows all become 14
$wsElements.Range("A2:D9").RowHeight = 14

# ---------------------------------------------------------------------------
# Sheet "Connections"
# ---------------------------------------------------------------------------
$wsConnections = $wb.Worksheets.Item("Connections")

# Column widths
$wsConnections.Columns.Item(1).ColumnWidth = 22.998697916666668
$wsConnections.Columns.Item(2).ColumnWidth = 23.330729166666668

# ---------------------------------------------------------------------------
# Sheet "Interactions"
# ---------------------------------------------------------------------------
$wsInteractions = $wb.Worksheets.Item("Interactions")

# Clear the bold/alt-font style from A2/B2 (keep C2's plain style) while
# preserving their values, then update the interaction-term values.
$styleSource = $wsInteractions.Range("C2")
$styleSource.Copy()
$wsInteractions.Range("A2").PasteSpecial(-4122)
$styleSource.Copy()
$wsInteractions.Range("B2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$wsInteractions.Range("B2").Value = "Perceived stress"
$wsInteractions.Range("C2").Value = "Depressive symptoms"

# ---------------------------------------------------------------------------
# Selections / active sheet: end on "Elements" so it is the active tab.
# ---------------------------------------------------------------------------
$wsInteractions.Activate()
$wsInteractions.Columns.Item(4).Select() | Out-Null

$wsConnections.Activate()
$wsConnections.Range("B10").Select() | Out-Null

$wsElements.Activate()
$wsElements.Range("B8").Select() | Out-Null
